$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Logs": append the new inbound-mail row (row 18)
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A18").Value = "Technische storing"
$logs.Range("B18").Value = "mailmind.test@zohomail.eu"
$logs.Range("C18").Value = "De website werkt niet goed. Is hier iets mis mee?"
$logs.Range("D18").Value = "IT / Technisch probleem"
$logs.Range("E18").Value = "Beste klant,`nBedankt voor uw bericht. Om uw probleem met de website beter te kunnen onderzoeken, hebben we wat meer informatie nodig. Kunt u specifiek aangeven welke problemen u ervaart bij het gebruik van de website? Krijgt u foutmeldingen te zien, of is er een specifieke pagina die niet laadt? Eventuele details die u kunt verstrekken, zoals de tijd waarop het probleem zich voordeed, zullen ons helpen om het probleem snel op te lossen.`nMet vriendelijke groet,`n[Je naam]`nE-mailassistent"
$logs.Range("F18").Value = "2025-06-24 20:10:34"
$logs.Range("G18").Value = "Ja"

# Undo the row-18 auto-height bump that Excel applies when a multi-line
# value lands in a row that previously had no explicit height.
$logs.Rows.Item(18).AutoFit()

# Extend the conditional-formatting ranges so they keep covering the
# whole table (D2:D17 -> D2:D18, G2:G17 -> G2:G18).
$logs.Range("D2:D17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D18"))
$logs.Range("G2:G17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G18"))

# ---------------------------------------------------------------------------
# Sheet "Dashboard": the category counts were refreshed from the logs, so
# "IT / Technisch probleem" rises to 3 occurrences and re-sorts into row 3,
# pushing the other categories down one row each.
# ---------------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "IT / Technisch probleem"
$dash.Range("B3").Value = 3

$dash.Range("A4").Value = "Sollicitatie / Vacature"
$dash.Range("B4").Value = 2

$dash.Range("A5").Value = "Factuur / Administratie"
$dash.Range("B5").Value = 2

$dash.Range("A6").Value = "Offerte / Prijsaanvraag"
$dash.Range("B6").Value = 2

$dash.Range("A7").Value = "Productinformatie"
$dash.Range("B7").Value = 2
